# BlancoValueObjectKtFieldStructure.xlsx
# 0.0.9: Adapt kotlin preferred type and annotation specifications.
#
# - Insert 3 new rows (typeKt, genericKt, annotationListKt) right after the
#   "default" row, shifting defaultKt/abstract/nullable/value/constArg and
#   everything below it down by 3 rows.
# - Fix up the "annotationList" description text (drop the "(excluding @)"
#   parenthetical, since that now only applies to the plain Java field).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("valueObject")

# --- Insert 3 blank rows right below row 32 ("default") -------------------
$ws.Rows("33:35").Insert()

# Copy formatting for the two plain "String"-typed rows (33,34) from the
# existing "type" row (29), and the List<String>-typed row (35, which needs
# the taller ht=45 style) from the existing "annotationList" row (31).
$ws.Range("A29:F29").Copy()
$ws.Range("A33:F34").PasteSpecial(-4122)

$ws.Range("A31:F31").Copy()
$ws.Range("A35:F35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 33: typeKt ---------------------------------------------------------
$ws.Range("A33").Formula = "=A32+1"
$ws.Range("B33").Value = "typeKt"
$ws.Range("C33").Value = "java.lang.String"
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = "Kotlin優先型名をパッケージ名のフル修飾付で指定します。"
$ws.Range("G33").Value = ""

# --- Row 34: genericKt ------------------------------------------------------
$ws.Range("A34").Formula = "=A33+1"
$ws.Range("B34").Value = "genericKt"
$ws.Range("C34").Value = "java.lang.String"
$ws.Range("D34").Value = ""
$ws.Range("E34").Value = "Kotlin優先型が期待する総称型の具体的な型名を指定します．"
$ws.Range("G34").Value = ""

# --- Row 35: annotationListKt -----------------------------------------------
$ws.Range("A35").Formula = "=A34+1"
$ws.Range("B35").Value = "annotationListKt"
$ws.Range("C35").Value = "java.util.List<java.lang.String>"
$ws.Range("D35").Value = "new java.util.ArrayList<java.lang.String>()"
$ws.Range("E35").Value = "Kotlin優先アノテーション文字列です"
$ws.Range("G35").Value = ""

# --- Fix the existing "annotationList" row's description -------------------
$ws.Range("E31").Value = "アノテーション文字列です"

$wb.Save()
